$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "65.801.84"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.54%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.675.26"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +0.82%  "
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "602.34"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.95%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "157.10"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.63%  "
$ws.Range("E7").Value = "  -0.02%  "
$ws.Range("E8").Value = "  +5.14%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.125"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +1.88%  "
$ws.Range("E10").Value = "  +1.07%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.401"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.36%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.155"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.25%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "29.49"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -1.37%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.0000198"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +1.08%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.155.63"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.77%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "65.603.06"
$ws.Range("D16").Style = "Normal"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.677.47"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +1.04%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "12.67"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.56%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.59"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +1.83%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "352.55"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -1.67%  "
$ws.Range("E22").Value = "  -0.12%  "
$ws.Range("E23").Value = "  +0.08%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.0000111"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +5.60%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "9.81"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +3.70%  "
$ws.Range("E26").Value = "  -5.23%  "
$ws.Range("E27").Value = "  +1.25%  "
$ws.Range("E28").Value = "  -1.52%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "8.14"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.13%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "545.06"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +4.00%  "
$ws.Range("E31").Value = "  +0.04%  "
$ws.Range("E32").Value = "  -1.68%  "
$ws.Range("E33").Value = "  -0.34%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "6.60"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +3.83%  "
$ws.Range("E35").Value = "  -1.51%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.424"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -1.88%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "20.42"
$ws.Range("D37").Style = "Normal"
$ws.Range("E38").Value = "  -0.02%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "159.12"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -2.24%  "
$ws.Range("E40").Value = "  -1.01%  "
$ws.Range("E41").Value = "  +0.02%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "42.79"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +2.19%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "165.14"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.30%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "4.10"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.71%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0614"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.34%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.31"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -1.27%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "23.26"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +1.29%  "
$ws.Range("E48").Value = "  -1.04%  "
$ws.Range("E49").Value = "  -1.23%  "
$ws.Range("E50").Value = "  +3.11%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "20.22"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +2.46%  "
